# "Add files via upload" / "a few changes to the table in the main table"
# Populate row 1 of Sheet1 with the Opportunities table data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 111111111
$ws.Range("B1").Value = "David"
$ws.Range("C1").Value = "Davidson"
$ws.Range("D1").Value = 546666666
$ws.Range("E1").Value = "New"
$ws.Range("F1").Value = 208063511

# G1 = 10/28/2018, stored as the Excel date serial number 43401
$ws.Range("G1").Value = 43401
$ws.Range("G1").NumberFormat = "mm-dd-yy"

$ws.Range("H1").Value = "don’t know what he wants in his life"

# Columns sized to fit their (now populated) content, as Excel's AutoFit would.
$ws.Columns.Item(1).ColumnWidth = 10.166666666666666
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 9.166666666666666
$ws.Columns.Item(7).ColumnWidth = 9.833333333333334
$ws.Columns.Item(8).ColumnWidth = 32.666666666666664

# Leave the selection where the author's session left it when saving.
$ws.Range("K7").Select()
